$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

# Lower "min_hours_per_week" (column I) from 30 to 25 for every employee row (2-58)
$ws.Range("I2:I58").Value = 25

# Remove the last three employees (Tina, Simon, Nicole) - rows 59, 60, 61
$ws.Range("A59:O61").EntireRow.Delete()

# Update the sheet view: scroll position and current selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J54").Select()
